$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2021" column (R) mirroring the formatting of the existing
# "2020" column (Q), then fill in the header year and the data value.

# Header cell R4 (year 2021), copy format from Q4 then set value
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Value = 2021

# Data cell R5 (value 42.9), copy format from Q5 then set value
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
$ws.Range("R5").Value = 42.9

# Update the active selection to match the new "last cell" reference
[void]$ws.Range("R9").Select()
